$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content changes -------------------------------------------------
# J2 previously held "Cfg"; it now holds "Spc" (a new shared string).
$ws.Range("J2").Value = "Spc"

# U2 was empty; it now holds "Cfg" (reusing the style used by the rest of
# row 2 in that block, as seen on T2 - yellow fill, large Consolas font).
$ws.Range("T2").Copy()
$ws.Range("U2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("U2").Value = "Cfg"

# AF6 was empty; it now holds "Const" (reusing the style used by AA3 -
# small Consolas font, yellow fill).
$ws.Range("AA3").Copy()
$ws.Range("AF6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AF6").Value = "Const"

$excel.CutCopyMode = 0

# --- View / selection changes ---------------------------------------------
# Scroll the sheet one column to the right and move the active selection
# from U9 to V9, as recorded in the saved workbook view.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("V9").Select()

$wb.Save()
